$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row additions (row 1), copying style from existing header cell E1
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

$ws.Range("E1").Copy()
$ws.Range("F1:H1").PasteSpecial(-4122)  # xlPasteFormats

# Re-set the values after pasting formats (paste formats should not affect value, but just to be safe)
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# Boolean FALSE values for rows 2-4
$ws.Range("F2:H4").Value = $false
